$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value()
    if ($v -ne $null) {
        $oa = $v.ToOADate()
        if ($oa -eq 46075) {
            $cell.Value = 46076
        }
    }
}
